$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B/C/D/E for reordered coin rows (36-45 block)
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '0.9012'
$ws.Range("E36").Value = '  -3.12%  '

$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.173.55'
$ws.Range("E37").Value = '  -0.23%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.646'
$ws.Range("E38").Value = '  -4.19%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").Value = '0.5657'
$ws.Range("E39").Value = '  -0.43%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '1.009'
$ws.Range("E41").Value = '  +0.34%  '

$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").Value = '2.561'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '100.12'
$ws.Range("E44").Value = '  -0.40%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8094'
$ws.Range("E45").Value = '  -3.92%  '

# Update D (and E where changed) for remaining rows
$ws.Range("D2").Value = '26.059.91'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '1.644.78'
$ws.Range("E3").Value = '  -0.97%  '

$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").Value = '215.91'
$ws.Range("E5").Value = '  -1.48%  '

$ws.Range("D6").Value = '0.5058'
$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("D7").Value = '1.009'
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("D8").Value = '0.06464'

$ws.Range("D9").Value = '0.2574'
$ws.Range("E9").Value = '  -0.76%  '

$ws.Range("D10").Value = '19.59'
$ws.Range("E10").Value = '  -2.18%  '

$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  -0.98%  '

$ws.Range("D12").Value = '1.642.87'
$ws.Range("E12").Value = '  -1.61%  '

$ws.Range("D13").Value = '4.253'
$ws.Range("E13").Value = '  -1.48%  '

$ws.Range("D14").Value = '1.867.55'
$ws.Range("E14").Value = '  -1.10%  '

$ws.Range("D15").Value = '0.5469'
$ws.Range("E15").Value = '  -1.70%  '

$ws.Range("D16").Value = '0.0₅7969'
$ws.Range("E16").Value = '  -1.37%  '

$ws.Range("D17").Value = '63.66'
$ws.Range("E17").Value = '  -1.27%  '

$ws.Range("D18").Value = '26.058.58'
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").Value = '1.009'
$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").Value = '205.34'
$ws.Range("E20").Value = '  -3.76%  '

$ws.Range("D21").Value = '4.310'
$ws.Range("E21").Value = '  -2.65%  '

$ws.Range("D22").Value = '10.04'
$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("D23").Value = '6.015'
$ws.Range("E23").Value = '  +1.13%  '

$ws.Range("D24").Value = '1.009'
$ws.Range("E24").Value = '  +0.35%  '

$ws.Range("D25").Value = '1.953'
$ws.Range("E25").Value = '  +11.03%  '

$ws.Range("D26").Value = '142.12'
$ws.Range("E26").Value = '  -1.58%  '

$ws.Range("D27").Value = '0.1156'
$ws.Range("E27").Value = '  -1.15%  '

$ws.Range("D28").Value = '15.79'
$ws.Range("E28").Value = '  -0.24%  '

$ws.Range("D29").Value = '6.751'
$ws.Range("E29").Value = '  -3.75%  '

$ws.Range("D30").Value = '0.05093'
$ws.Range("E30").Value = '  -4.24%  '

$ws.Range("D31").Value = '1.245'

$ws.Range("D32").Value = '3.269'
$ws.Range("E32").Value = '  -3.16%  '

$ws.Range("D33").Value = '3.194'
$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").Value = '1.548'
$ws.Range("E34").Value = '  -2.11%  '

$ws.Range("D35").Value = '2.359'
$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("D40").Value = '0.01577'
$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D43").Value = '5.692'
$ws.Range("E43").Value = '  +0.38%  '

$ws.Range("D46").Value = '1.779.21'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("D47").Value = '0.0₈112'
$ws.Range("E47").Value = '  -1.25%  '

$ws.Range("D48").Value = '0.4537'
$ws.Range("E48").Value = '  +0.23%  '

$ws.Range("D49").Value = '1.010'
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("D50").Value = '55.26'
$ws.Range("E50").Value = '  -1.61%  '

$ws.Range("D51").Value = '0.05040'
$ws.Range("E51").Value = '  -0.82%  '
